$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.34488412926216
$ws.Range("C2").Value = 10.53754694196525
$ws.Range("D2").Value = 5.998283038567663
$ws.Range("E2").Value = 13.088582084101
$ws.Range("G2").Value = 34.60046559357446
$ws.Range("H2").Value = 15.81609794807873
$ws.Range("K2").Value = 8.743748939633958
$ws.Range("L2").Value = 9.916044921239543
$ws.Range("M2").Value = 14.06412561640813
$ws.Range("N2").Value = 20.20494549051299
$ws.Range("O2").Value = 24.84269747074645
$ws.Range("B3").Value = 12.13014762780094
$ws.Range("C3").Value = 10.53687347929094
$ws.Range("D3").Value = 5.881568747071912
$ws.Range("E3").Value = 13.11777255556099
$ws.Range("G3").Value = 34.65023836534969
$ws.Range("H3").Value = 15.85780999390773
$ws.Range("K3").Value = 8.577283378249001
$ws.Range("L3").Value = 9.922994243321313
$ws.Range("M3").Value = 14.03464532889443
$ws.Range("N3").Value = 20.26382100262406
$ws.Range("O3").Value = 24.90575112223637
$ws.Range("B4").Value = 11.99873679461833
$ws.Range("C4").Value = 10.53675326343288
$ws.Range("D4").Value = 5.810450533656493
$ws.Range("E4").Value = 13.13742495078687
$ws.Range("G4").Value = 34.68970678908941
$ws.Range("H4").Value = 15.8856479508427
$ws.Range("K4").Value = 8.474846217240705
$ws.Range("L4").Value = 9.92857355745571
$ws.Range("M4").Value = 14.01858958889535
$ws.Range("N4").Value = 20.30167064027578
$ws.Range("O4").Value = 24.9490146918185
$ws.Range("B5").Value = 11.94536795603046
$ws.Range("C5").Value = 10.53677843628117
$ws.Range("D5").Value = 5.781651152467171
$ws.Range("E5").Value = 13.14586860470865
$ws.Range("G5").Value = 34.70802609415936
$ws.Range("H5").Value = 15.89755217712374
$ws.Range("K5").Value = 8.433101588323026
$ws.Range("L5").Value = 9.931177781269209
$ws.Range("M5").Value = 14.01256600316368
$ws.Range("N5").Value = 20.31752334209003
$ws.Range("O5").Value = 24.96778746740788
$ws.Range("B6").Value = 11.93651921758548
$ws.Range("C6").Value = 10.53678710738985
$ws.Range("D6").Value = 5.776881340208949
$ws.Range("E6").Value = 13.14729695682007
$ws.Range("G6").Value = 34.71120288168336
$ws.Range("H6").Value = 15.89956269335421
$ws.Range("K6").Value = 8.426171553186707
$ws.Range("L6").Value = 9.931630194636284
$ws.Range("M6").Value = 14.0115972909392
$ws.Range("N6").Value = 20.32018159778107
$ws.Range("O6").Value = 24.97097364059271
$ws.Range("B7").Value = 11.99801621075102
$ws.Range("C7").Value = 10.5367533020998
$ws.Range("D7").Value = 5.810061339170534
$ws.Range("E7").Value = 13.13753706259834
$ws.Range("G7").Value = 34.68994480417011
$ws.Range("H7").Value = 15.88580622767606
$ws.Range("K7").Value = 8.474283161963902
$ws.Range("L7").Value = 9.92860733957032
$ws.Range("M7").Value = 14.01850624397643
$ws.Range("N7").Value = 20.30188269806956
$ws.Range("O7").Value = 24.94926324369416
$ws.Range("B8").Value = 12.27079493419088
$ws.Range("C8").Value = 10.537254105035
$ws.Range("D8").Value = 5.957953635132041
$ws.Range("E8").Value = 13.09828821519575
$ws.Range("G8").Value = 34.61577610937546
$ws.Range("H8").Value = 15.83001821852253
$ws.Range("K8").Value = 8.686432070498231
$ws.Range("L8").Value = 9.918169032307656
$ws.Range("M8").Value = 14.05353939625256
$ws.Range("N8").Value = 20.22489371009615
$ws.Range("O8").Value = 24.86349366055815
$ws.Range("B9").Value = 12.80608099710406
$ws.Range("C9").Value = 10.54054573215861
$ws.Range("D9").Value = 6.250353991795447
$ws.Range("E9").Value = 13.03502971112377
$ws.Range("G9").Value = 34.54116645714353
$ws.Range("H9").Value = 15.73827884621939
$ws.Range("K9").Value = 9.098251494392358
$ws.Range("L9").Value = 9.908085026228994
$ws.Range("M9").Value = 14.13823541768264
$ws.Range("N9").Value = 20.08734883865613
$ws.Range("O9").Value = 24.7314398608239
$ws.Range("B10").Value = 13.19552713119145
$ws.Range("C10").Value = 10.54434908352278
$ws.Range("D10").Value = 6.464112385317796
$ws.Range("E10").Value = 12.99689254601012
$ws.Range("G10").Value = 34.52969125275679
$ws.Range("H10").Value = 15.681635625464
$ws.Range("K10").Value = 9.395155838884801
$ws.Range("L10").Value = 9.906967542477824
$ws.Range("M10").Value = 14.20988138893889
$ws.Range("N10").Value = 19.99440068424357
$ws.Range("O10").Value = 24.65651657997451
$ws.Range("B11").Value = 13.37104991791745
$ws.Range("C11").Value = 10.54637516432672
$ws.Range("D11").Value = 6.560620749481568
$ws.Range("E11").Value = 12.98134952291793
$ws.Range("G11").Value = 34.53389476749058
$ws.Range("H11").Value = 15.65820110362894
$ws.Range("K11").Value = 9.528384198229874
$ws.Range("L11").Value = 9.907815333796805
$ws.Range("M11").Value = 14.24444256962732
$ws.Range("N11").Value = 19.95385894019554
$ws.Range("O11").Value = 24.62724162828778
$ws.Range("B12").Value = 13.43721600239276
$ws.Range("C12").Value = 10.54718450948217
$ws.Range("D12").Value = 6.597021379967178
$ws.Range("E12").Value = 12.9757230917286
$ws.Range("G12").Value = 34.53684081754089
$ws.Range("H12").Value = 15.64966235783011
$ws.Range("K12").Value = 9.578522823244107
$ws.Range("L12").Value = 9.90833039965413
$ws.Range("M12").Value = 14.25780599139667
$ws.Range("N12").Value = 19.93875588733131
$ws.Range("O12").Value = 24.616848076682
$ws.Range("B13").Value = 13.42298032233226
$ws.Range("C13").Value = 10.54700833661174
$ws.Range("D13").Value = 6.589188912376859
$ws.Range("E13").Value = 12.97692331349145
$ws.Range("G13").Value = 34.53614612657592
$ws.Range("H13").Value = 15.65148641348296
$ws.Range("K13").Value = 9.567739190422099
$ws.Range("L13").Value = 9.90821085817872
$ws.Range("M13").Value = 14.25491579032172
$ws.Range("N13").Value = 19.94199753308894
$ws.Range("O13").Value = 24.61905571458359
$ws.Range("B14").Value = 13.37649982683857
$ws.Range("C14").Value = 10.54644090817976
$ws.Range("D14").Value = 6.563618577488317
$ws.Range("E14").Value = 12.98088143653828
$ws.Range("G14").Value = 34.53411000795366
$ws.Range("H14").Value = 15.65749189430787
$ws.Range("K14").Value = 9.532515653571064
$ws.Range("L14").Value = 9.907853825931008
$ws.Range("M14").Value = 14.24553650079409
$ws.Range("N14").Value = 19.95261141438163
$ws.Range("O14").Value = 24.62637266673235
$ws.Range("B15").Value = 13.34798816046247
$ws.Range("C15").Value = 10.54609881302486
$ws.Range("D15").Value = 6.547935981915666
$ws.Range("E15").Value = 12.98333967035587
$ws.Range("G15").Value = 34.53303914866141
$ws.Range("H15").Value = 15.66121410187624
$ws.Range("K15").Value = 9.510898196859301
$ws.Range("L15").Value = 9.907660370672595
$ws.Range("M15").Value = 14.23982711735461
$ws.Range("N15").Value = 19.95914514708682
$ws.Range("O15").Value = 24.63094468362936
$ws.Range("B16").Value = 13.18401771053808
$ws.Range("C16").Value = 10.54422259143792
$ws.Range("D16").Value = 6.457787166399072
$ws.Range("E16").Value = 12.99794463119661
$ws.Range("G16").Value = 34.52960612461272
$ws.Range("H16").Value = 15.6832140638076
$ws.Range("K16").Value = 9.386407945690308
$ws.Range("L16").Value = 9.90693934383868
$ws.Range("M16").Value = 14.20766174070404
$ws.Range("N16").Value = 19.9970851180006
$ws.Range("O16").Value = 24.6585265925678
$ws.Range("B17").Value = 13.08296032592107
$ws.Range("C17").Value = 10.54314706965024
$ws.Range("D17").Value = 6.402267452943692
$ws.Range("E17").Value = 13.00736659968514
$ws.Range("G17").Value = 34.52991345130038
$ws.Range("H17").Value = 15.69730777343008
$ws.Range("K17").Value = 9.309532077444509
$ws.Range("L17").Value = 9.906843680497019
$ws.Range("M17").Value = 14.1884285597193
$ws.Range("N17").Value = 20.02080514227212
$ws.Range("O17").Value = 24.6766793008719
$ws.Range("B18").Value = 13.02468429331654
$ws.Range("C18").Value = 10.54255633764432
$ws.Range("D18").Value = 6.370267390066005
$ws.Range("E18").Value = 13.01295585264931
$ws.Range("G18").Value = 34.53097733528789
$ws.Range("H18").Value = 15.70563367690843
$ws.Range("K18").Value = 9.265145113054192
$ws.Range("L18").Value = 9.906916367665618
$ws.Range("M18").Value = 14.17755193819195
$ws.Range("N18").Value = 20.03461219408073
$ws.Range("O18").Value = 24.68757277081934
$ws.Range("B19").Value = 13.00492921370148
$ws.Range("C19").Value = 10.54236112765217
$ws.Range("D19").Value = 6.359422529580098
$ws.Range("E19").Value = 13.0148774814785
$ws.Range("G19").Value = 34.53148991890036
$ws.Range("H19").Value = 15.70849039819935
$ws.Range("K19").Value = 9.250088765914152
$ws.Range("L19").Value = 9.906962941407562
$ws.Range("M19").Value = 14.17390142868497
$ws.Range("N19").Value = 20.03931521205175
$ws.Range("O19").Value = 24.69133880375216
$ws.Range("B20").Value = 13.09373411386334
$ws.Range("C20").Value = 10.54325867851879
$ws.Range("D20").Value = 6.40818480852356
$ws.Range("E20").Value = 13.00634602421044
$ws.Range("G20").Value = 34.5297889250268
$ws.Range("H20").Value = 15.69578474915869
$ws.Range("K20").Value = 9.317733596952001
$ws.Range("L20").Value = 9.90684065339477
$ws.Range("M20").Value = 14.19045678738396
$ws.Range("N20").Value = 20.01826314799191
$ws.Range("O20").Value = 24.67470007414223
$ws.Range("B21").Value = 13.39016093531832
$ws.Range("C21").Value = 10.54660643627493
$ws.Range("D21").Value = 6.571133446311695
$ws.Range("E21").Value = 12.97971180294223
$ws.Range("G21").Value = 34.53467132302404
$ws.Range("H21").Value = 15.65571883684156
$ws.Range("K21").Value = 9.542870499927405
$ws.Range("L21").Value = 9.907953437393328
$ws.Range("M21").Value = 14.24828399863535
$ws.Range("N21").Value = 19.94948710499749
$ws.Range("O21").Value = 24.62420470735781
$ws.Range("B22").Value = 13.58211307040566
$ws.Range("C22").Value = 10.54903964565136
$ws.Range("D22").Value = 6.676767758955807
$ws.Range("E22").Value = 12.96381648128398
$ws.Range("G22").Value = 34.54575480088399
$ws.Range("H22").Value = 15.63148837671804
$ws.Range("K22").Value = 9.688171148367978
$ws.Range("L22").Value = 9.909811124445168
$ws.Range("M22").Value = 14.28768170180059
$ws.Range("N22").Value = 19.90599015830113
$ws.Range("O22").Value = 24.59523824422006
$ws.Range("B23").Value = 13.47984811950628
$ws.Range("C23").Value = 10.54771870139369
$ws.Range("D23").Value = 6.620480203925763
$ws.Range("E23").Value = 12.97216189662684
$ws.Range("G23").Value = 34.5391177356871
$ws.Range("H23").Value = 15.64424177292454
$ws.Range("K23").Value = 9.610804962977291
$ws.Range("L23").Value = 9.908716554125373
$ws.Range("M23").Value = 14.26651010571888
$ws.Range("N23").Value = 19.92907278565871
$ws.Range("O23").Value = 24.61032872916733
$ws.Range("B24").Value = 13.08886383002906
$ws.Range("C24").Value = 10.54320813415093
$ws.Range("D24").Value = 6.405509821446219
$ws.Range("E24").Value = 13.00680688924075
$ws.Range("G24").Value = 34.52984245974508
$ws.Range("H24").Value = 15.69647261303752
$ws.Range("K24").Value = 9.314026278101954
$ws.Range("L24").Value = 9.906841624151161
$ws.Range("M24").Value = 14.18953926196411
$ws.Range("N24").Value = 20.01941185368064
$ws.Range("O24").Value = 24.67559345834625
$ws.Range("B25").Value = 12.66165892499734
$ws.Range("C25").Value = 10.53941036614441
$ws.Range("D25").Value = 6.171262501443554
$ws.Range("E25").Value = 13.050677042667
$ws.Range("G25").Value = 34.55374669309602
$ws.Range("H25").Value = 15.76120694125001
$ws.Range("K25").Value = 8.987621305450974
$ws.Range("L25").Value = 9.909705023933956
$ws.Range("M25").Value = 14.11364324033807
$ws.Range("N25").Value = 20.12312901859702
$ws.Range("O25").Value = 24.76328769911588

Write-Output "Applied 380 kV case update (264 cells)"
